$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), reusing the same header style as the
# existing sum column (G1) so no new style entry is introduced.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# New data cell for the Save column (H2)
$ws.Range("H2").Value = 0
